$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Price (column D) and Volume(1h) (column E) values for the
# cryptos list, as refreshed by the scheduled GitHub Actions run.
$updates = @(
    @{Row=2; D="68.781.00"; E="  -0.69%  "},
    @{Row=3; D="3.935.44"; E="  +3.75%  "},
    @{Row=4; D="1.00"; E="  -0.06%  "},
    @{Row=5; D="604.04"; E="  +0.04%  "},
    @{Row=6; D="166.20"; E="  +0.72%  "},
    @{Row=7; D="3.932.27"; E="  +3.81%  "},
    @{Row=8; D=$null; E="  -0.07%  "},
    @{Row=9; D="0.530"; E="  -1.47%  "},
    @{Row=10; D="0.166"; E="  -2.76%  "},
    @{Row=11; D=$null; E="  +0.50%  "},
    @{Row=12; D="0.461"; E="  +0.10%  "},
    @{Row=13; D=$null; E="  +0.88%  "},
    @{Row=14; D="37.20"; E="  -0.36%  "},
    @{Row=15; D="4.593.72"; E="  +3.72%  "},
    @{Row=16; D="3.890.49"; E="  +2.37%  "},
    @{Row=17; D="68.928.09"; E="  -0.60%  "},
    @{Row=18; D=$null; E="  -0.12%  "},
    @{Row=19; D="17.12"; E="  -2.00%  "},
    @{Row=20; D=$null; E="  -1.36%  "},
    @{Row=21; D="11.18"; E="  -0.28%  "},
    @{Row=22; D="485.58"; E="  -1.56%  "},
    @{Row=23; D=$null; E="  +14.01%  "},
    @{Row=24; D="0.723"; E="  -0.13%  "},
    @{Row=25; D="84.81"; E="  -0.03%  "},
    @{Row=26; D="2.24"; E="  -1.28%  "},
    @{Row=27; D="12.02"; E="  -2.11%  "},
    @{Row=28; D="10.10"; E="  -0.16%  "},
    @{Row=29; D=$null; E="  -0.05%  "},
    @{Row=30; D="2.93"; E="  -1.86%  "},
    @{Row=31; D="4.088.45"; E="  +3.67%  "},
    @{Row=32; D="2.39"; E="  -0.46%  "},
    @{Row=33; D="32.22"; E="  +0.53%  "},
    @{Row=34; D="7.80"; E="  -3.94%  "},
    @{Row=35; D="3.884.95"; E="  +3.84%  "},
    @{Row=36; D="0.106"; E="  -0.66%  "},
    @{Row=37; D=$null; E="  +2.16%  "},
    @{Row=38; D=$null; E="  -1.17%  "},
    @{Row=39; D="5.91"; E="  -0.58%  "},
    @{Row=40; D="3.18"; E="  +4.12%  "},
    @{Row=41; D=$null; E="  +0.05%  "},
    @{Row=42; D="0.317"; E="  -2.02%  "},
    @{Row=43; D="433.84"; E="  +1.88%  "},
    @{Row=44; D="48.49"; E="  +0.10%  "},
    @{Row=45; D=$null; E="  -0.33%  "},
    @{Row=46; D="8.52"; E="  +1.15%  "},
    @{Row=47; D=$null; E="  +0.00%  "},
    @{Row=48; D="26.76"; E="  +10.24%  "},
    @{Row=49; D="2.838.35"; E="  +0.79%  "},
    @{Row=50; D="141.75"; E="  -0.39%  "},
    @{Row=51; D="0.000266"; E="  +18.72%  "}
)


foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Force text formatting first so Excel doesn't silently reinterpret
        # strings like "1.00" or "604.04" as numbers (which would drop the
        # trailing zeros / thousands-dot formatting used in this sheet).
        $dCell = $ws.Cells.Item($u.Row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        # Revert the cell style so no stray formatting is left behind.
        $dCell.Style = "Normal"
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
